$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "0.163", "85.09") are preserved as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '69.469.35'
$ws.Range('E2').Value = '  +2.90%  '
$ws.Range('D3').Value = '3.813.00'
$ws.Range('E3').Value = '  +1.36%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '599.82'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').Value = '171.91'
$ws.Range('E6').Value = '  +1.11%  '
$ws.Range('D7').Value = '3.811.35'
$ws.Range('E7').Value = '  +1.32%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -0.15%  '
$ws.Range('D10').Value = '0.163'
$ws.Range('E10').Value = '  -1.52%  '
$ws.Range('D11').Value = '6.55'
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').Value = '0.0000266'
$ws.Range('E13').Value = '  -4.45%  '
$ws.Range('D14').Value = '36.95'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('D15').Value = '4.451.12'
$ws.Range('E15').Value = '  +1.47%  '
$ws.Range('D16').Value = '3.801.79'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').Value = '69.397.06'
$ws.Range('E17').Value = '  +2.88%  '
$ws.Range('E18').Value = '  -2.62%  '
$ws.Range('D19').Value = '7.12'
$ws.Range('E19').Value = '  -1.47%  '
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('E21').Value = '  +5.56%  '
$ws.Range('D22').Value = '473.59'
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('E23').Value = '  -1.44%  '
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').Value = '0.0000149'
$ws.Range('E24').Value = '  +0.78%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '85.09'
$ws.Range('E25').Value = '  +1.46%  '
$ws.Range('D26').Value = '2.25'
$ws.Range('E26').Value = '  +0.94%  '
$ws.Range('D27').Value = '12.27'
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('D28').Value = '10.28'
$ws.Range('E28').Value = '  -1.08%  '
$ws.Range('D30').Value = '3.960.93'
$ws.Range('E30').Value = '  +1.29%  '
$ws.Range('E31').Value = '  -2.50%  '
$ws.Range('D32').Value = '7.50'
$ws.Range('E32').Value = '  -2.85%  '
$ws.Range('E33').Value = '  +0.48%  '
$ws.Range('D34').Value = '30.40'
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('D35').Value = '9.43'
$ws.Range('E35').Value = '  +2.92%  '
$ws.Range('D37').Value = '3.765.67'
$ws.Range('E37').Value = '  +1.17%  '
$ws.Range('E38').Value = '  -2.34%  '
$ws.Range('D39').Value = '3.56'
$ws.Range('E39').Value = '  -7.18%  '
$ws.Range('E40').Value = '  +1.96%  '
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = '2.00'
$ws.Range('E46').Value = '  +2.34%  '
$ws.Range('D47').Value = '43.96'
$ws.Range('E47').Value = '  +12.77%  '
$ws.Range('D48').Value = '8.67'
$ws.Range('E48').Value = '  -1.09%  '
$ws.Range('E49').Value = '  +1.25%  '
$ws.Range('D50').Value = '404.33'
$ws.Range('D51').Value = '145.86'
$ws.Range('E51').Value = '  +3.04%  '
